# Auto-generated edit script: update cryptos Price (D) / Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.273.28'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '2.444.82'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.69'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.60'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.67%  '
$ws.Range("D9").Value = '2.443.25'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  -3.02%  '
$ws.Range("E11").Value = '  +1.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.21'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.43'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("E15").Value = '  -3.58%  '
$ws.Range("D16").Value = '2.878.16'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '62.114.19'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '2.439.23'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.12'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.66'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.97'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.89%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.86'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.38'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '619.15'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").Value = '2.564.22'
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").Value = '0.0₃0955'
$ws.Range("E29").Value = '  -7.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  -4.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.01'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.38%  '
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.92'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.83%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.44'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.375'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '151.35'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.34'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.24'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.77'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.56'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.60%  '
$ws.Range("E45").Value = '  -8.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '143.42'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.63'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0525'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.53'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -7.84%  '
$ws.Range("E51").Value = '  -1.25%  '
